$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("A11").NumberFormat = "m/d/yyyy"
$ws.Range("A11").Value = 42795
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = "Updating SQL tables and Report"

# Row 12
$ws.Range("A12").NumberFormat = "m/d/yyyy"
$ws.Range("A12").Value = 42796
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "Updating Rating System Schema"

# Update selection to C13
$ws.Range("C13").Select()
